# "verificação de desconexão com bluetooth melhorada"
#
# Slide 1 ("Imagem 4" / shape id 5, the last Bluetooth-icon picture in the
# row) is duplicated. The duplicate is moved further down the slide and is
# given a solid red fill so the new (disconnected) state icon stands out,
# while keeping the same picture crop/size as the original.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the source picture to clone ("Imagem 4", the right-most icon).
$src = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Imagem 4") {
        $src = $candidate
    }
}

$newPic = $src.Duplicate().Item(1)
$newPic.Name = "Imagem 11"

# Keep the original size, move the duplicate down below the icon row.
$newPic.Left = 9117195 / 12700
$newPic.Top = 4893017 / 12700
$newPic.Width = 2560320 / 12700
$newPic.Height = 2156143 / 12700

# Flag the new icon red (disconnected bluetooth state).
$newPic.Fill.Visible = $true
$newPic.Fill.ForeColor.RGB = 255
